# Update the two "row 2" sample/test data rows in the invoice verification
# workbook with new test invoice values (per commit: config.properties /
# testng.xml test-data refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet "Memo_Verification_details" (PO verification row) ---
$wsVerification = $wb.Worksheets.Item("Memo_Verification_details")

# Invoice_number
$wsVerification.Range("B2").Value = "TESTINV97479"

# PO Number (numeric-looking text -> force text storage so it doesn't
# get coerced into a numeric cell)
$wsVerification.Range("C2").NumberFormat = "@"
$wsVerification.Range("C2").Value = "4500000891"

# Revised Tax Code
$wsVerification.Range("K2").Value = "KG"

# --- Sheet "Memo_invoice_Details" (invoice detail row) ---
$wsInvoice = $wb.Worksheets.Item("Memo_invoice_Details")

# Invoice_number
$wsInvoice.Range("B2").Value = "TESTINV97479"

# Invoice_Date (force text so the literal date string is preserved
# instead of being converted to a date serial number)
$wsInvoice.Range("C2").NumberFormat = "@"
$wsInvoice.Range("C2").Value = "2024-03-27"

# IGST (force text so the literal "0.18" is preserved instead of being
# converted to a floating point number)
$wsInvoice.Range("I2").NumberFormat = "@"
$wsInvoice.Range("I2").Value = "0.18"

# Total_Invoice_Amount (force text so the literal "1.18" is preserved)
$wsInvoice.Range("O2").NumberFormat = "@"
$wsInvoice.Range("O2").Value = "1.18"
